# JournalDeBord.xlsx update: "updated logbook + backlog"
#
# This script updates the single worksheet (Feuil1 / ActiveSheet) to:
#  - shorten the "N° question" header to "N°"
#  - add a handful of new Question/Réponse rows to the backlog table (cols F/G)
#  - update the text of the "pour sécuriser le dossier ..." backlog entry
#  - adjust two existing time entries and append two new logbook entries
#
# New unique strings are introduced in the exact order below so that they are
# appended to the shared string table in the same order as in the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new logbook rows appended at the bottom of the left-hand journal table ---
$ws.Range("D38").Value = "Mot de passe"
$ws.Range("C38").Value = "Rdv avec M Borys F"
$ws.Range("B38").Value = 0.5625
$ws.Range("B38").NumberFormat = "h:mm"

$ws.Range("B37").Value = 0.54166666666666663
$ws.Range("B37").NumberFormat = "h:mm"
$ws.Range("C37").Value = "Documentation"

# --- backlog entry text update ---
$ws.Range("D26").Value = "pour sécuriser le dossier les information de la DB"

# --- shorten the Q&A table header ---
$ws.Range("E1").Value = "N°"

# --- new Question / Réponse rows in the Q&A table (columns F/G) ---
$ws.Range("F3").Value = "Une utlisateur connecter rajoute en plus a la page d'aceuille ses favoris ou on remplace tout par les favoris ?"
$ws.Range("G3").Value = "On ne fais que rajouter un champs favoris"

$ws.Range("F4").Value = "Français ou en anglais le code ?"
$ws.Range("G4").Value = "Pas d’importance mais ne pas mélanger les deux."

$ws.Range("F5").Value = "Comment faire avec GANTT pour les vendredis après-midi vu que GANTT ne peut pas faire moins qu’une journée."

$ws.Range("F6").Value = "MCD "

# --- adjust the time of two existing logbook entries ---
$ws.Range("B35").Value = 0.41666666666666669
$ws.Range("B36").Value = 0.44097222222222227

# --- column widths for the now-visible Q&A / backlog columns ---
$ws.Columns.Item(4).ColumnWidth = 38.5
$ws.Columns.Item(5).ColumnWidth = 2.1666666666666665
$ws.Columns.Item(6).ColumnWidth = 89.5
$ws.Columns.Item(7).ColumnWidth = 38.666666666666664

# --- update the view to match the author's final selection/scroll position ---
$ws.Range("F6").Select()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
